$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column D (Price) temporarily to Text format so numeric-looking
# strings (e.g. "1.003", "10.00") are preserved as literal text instead
# of being auto-converted to real numbers by Excel, matching the source
# workbook where every data cell is an inline string.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.470.50'
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").Value = '1.830.98'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("E4").Value = '  -0.80%  '
$ws.Range("D5").Value = '331.62'
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("D7").Value = '0.4579'
$ws.Range("E7").Value = '  -2.69%  '
$ws.Range("D8").Value = '0.3808'
$ws.Range("E8").Value = '  -2.94%  '
$ws.Range("D9").Value = '46.44'
$ws.Range("E9").Value = '  +2.96%  '
$ws.Range("D10").Value = '0.07905'
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("D11").Value = '0.9699'
$ws.Range("E11").Value = '  -3.44%  '
$ws.Range("D12").Value = '21.08'
$ws.Range("E12").Value = '  -3.64%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.890'
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.826.58'
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("D15").Value = '7.074'
$ws.Range("E15").Value = '  -2.63%  '
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").Value = '89.36'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '0.06613'
$ws.Range("E18").Value = '  -1.91%  '
$ws.Range("D19").Value = '0.00001027'
$ws.Range("E19").Value = '  -1.66%  '
$ws.Range("D20").Value = '17.14'
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  -0.61%  '
$ws.Range("D22").Value = '27.457.15'
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("D23").Value = '5.341'
$ws.Range("E23").Value = '  -2.37%  '
$ws.Range("D24").Value = '10.81'
$ws.Range("E24").Value = '  -1.03%  '
$ws.Range("D25").Value = '2.296'
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("D26").Value = '2.044.28'
$ws.Range("E26").Value = '  -2.06%  '
$ws.Range("D27").Value = '155.43'
$ws.Range("E27").Value = '  -2.49%  '
$ws.Range("D28").Value = '19.41'
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").Value = '2.068'
$ws.Range("E29").Value = '  -4.32%  '
$ws.Range("D30").Value = '5.294'
$ws.Range("E30").Value = '  -3.02%  '
$ws.Range("D31").Value = '118.57'
$ws.Range("E31").Value = '  -2.75%  '
$ws.Range("D32").Value = '0.9414'
$ws.Range("E32").Value = '  -4.35%  '
$ws.Range("D33").Value = '0.09304'
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("D34").Value = '3.579'
$ws.Range("E34").Value = '  -0.88%  '
$ws.Range("E35").Value = '  -1.32%  '
$ws.Range("D36").Value = '1.334'
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("D37").Value = '0.05932'
$ws.Range("E37").Value = '  -2.19%  '
$ws.Range("D38").Value = '0.02187'
$ws.Range("E38").Value = '  -2.24%  '
$ws.Range("D39").Value = '8.088'
$ws.Range("E39").Value = '  -3.13%  '
$ws.Range("D40").Value = '1.151'
$ws.Range("E40").Value = '  -3.84%  '
$ws.Range("D41").Value = '0.5781'
$ws.Range("E41").Value = '  -3.43%  '
$ws.Range("D42").Value = '0.1829'
$ws.Range("E42").Value = '  -3.15%  '
$ws.Range("D43").Value = '10.00'
$ws.Range("E43").Value = '  -3.02%  '
$ws.Range("D44").Value = '1.277'
$ws.Range("E44").Value = '  +2.39%  '
$ws.Range("D45").Value = '0.5459'
$ws.Range("E45").Value = '  -3.53%  '
$ws.Range("D46").Value = '11.93'
$ws.Range("E46").Value = '  -3.04%  '
$ws.Range("D47").Value = '1.872'
$ws.Range("E47").Value = '  -2.82%  '
$ws.Range("D48").Value = '110.84'
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("D49").Value = '0.06608'
$ws.Range("E49").Value = '  -2.25%  '
$ws.Range("D50").Value = '1.003'
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("D51").Value = '1.045'
$ws.Range("E51").Value = '  -1.22%  '

# Restore the default (General) formatting so no stray cell style is
# introduced versus the original workbook.
$priceRange.ClearFormats()

